# Commit: Thu, May 14, 2020 11:05:54 AM
#
# 1) The table on slide 16 switches to a different built-in table style
#    (tableStyleId {61B2A5D6-29F2-4FED-885C-087BC2CA1364} ->
#     {E26D1D02-2FCE-42FA-AAD0-6626B2EC87F9}).
# 2) The deck's theme (ppt/theme/theme1.xml, the SlideMaster's theme that
#    drives the whole presentation's look) is switched from the custom
#    "Integral" palette to the stock Office palette ("Office Theme").

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 16 -------------------------------------
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{E26D1D02-2FCE-42FA-AAD0-6626B2EC87F9}", $true)
    }
}

# --- 2. Swap the presentation theme colors to the default Office palette ---
$officeColors = @(
    0,          # dk1     000000
    16777215,   # lt1     FFFFFF
    6968388,    # dk2     44546A
    15132391,   # lt2     E7E6E6
    13998939,   # accent1 5B9BD5
    3243501,    # accent2 ED7D31
    10855845,   # accent3 A5A5A5
    49407,      # accent4 FFC000
    12874308,   # accent5 4472C4
    4697456,    # accent6 70AD47
    12673797,   # hlink   0563C1
    7491477     # folHlink 954F72
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
